$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary row 12: average of column J (|S*|/n) ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# --- New summary rows 14-17: labelled aggregate statistics ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Formatting for the new B14:B17 values: bold, size 12, vertically centered ---
# Build the style on the first cell, then fan it out via a format-only paste so
# the style table doesn't accumulate unused transitional cellXfs entries.
$firstLabel = $ws.Range("B14")
$firstLabel.Font.Bold = $true
$firstLabel.Font.Size = 12
$firstLabel.Font.Charset = 161
$firstLabel.VerticalAlignment = -4108

$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15:B17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row heights for the labelled rows (visually taller for the bigger font)
$ws.Range("A14:B17").RowHeight = 15.6

# --- Page setup (paper size + orientation) as recorded by the author's Excel ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection left on the newly added summary block, matching the saved view ---
$ws.Range("A14:B17").Select() | Out-Null
